# CT graphed profiles 1/15 script and output
# Updates the Cond_temp_logger_CalLaunch workbook:
#  - corrects the "Log Start" time (col E) for the two existing 1/15 rows
#    that used probe 318 (rows 20-21)
#  - appends a new logged row (row 24) for probe 319 on 1/15, including a
#    new Notes comment
#  - moves the active selection to reflect where the user ended up after
#    entering the new row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix the Log Start time recorded for the two probe-318 rows on 1/15 ---
$ws.Range("E20").Value = 0.41530092592592593
$ws.Range("E21").Value = 0.41530092592592593

# --- append the new row for probe 319 logged on 1/15 ---
$ws.Range("A24").Value = 44211
$ws.Range("A24").NumberFormat = $ws.Range("A23").NumberFormat

$ws.Range("B24").Value = 319

$ws.Range("C24").Value = 0.36440972222222223
$ws.Range("C24").NumberFormat = $ws.Range("C23").NumberFormat

$ws.Range("D24").Value = 0.36660879629629628
$ws.Range("D24").NumberFormat = $ws.Range("D23").NumberFormat

$ws.Range("E24").Value = 0.3669560185185185
$ws.Range("E24").NumberFormat = $ws.Range("E23").NumberFormat

$ws.Range("F24").Value = 0.61880787037037044
$ws.Range("F24").NumberFormat = $ws.Range("F23").NumberFormat

$ws.Range("I24").Value = 50000
$ws.Range("I24").NumberFormat = $ws.Range("I23").NumberFormat

$ws.Range("J24").Value = "kept in calibration solution all day"

# --- reflect the new cursor position / selection left after data entry ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 2
$ws.Range("D25").Select()
